$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "26.678.74"
$ws.Cells.Item(2, 5).Value = "  -0.05%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.596.89"
$ws.Cells.Item(3, 5).Value = "  -0.14%  "

$ws.Cells.Item(4, 5).Value = "  +0.32%  "

Set-TextValue $ws.Cells.Item(5, 4) "211.32"
$ws.Cells.Item(5, 5).Value = "  +0.07%  "

Set-TextValue $ws.Cells.Item(6, 4) "0.512"
$ws.Cells.Item(6, 5).Value = "  -0.30%  "

$ws.Cells.Item(7, 5).Value = "  +0.32%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.0617"
$ws.Cells.Item(8, 5).Value = "  +0.01%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.246"
$ws.Cells.Item(9, 5).Value = "  +0.17%  "

Set-TextValue $ws.Cells.Item(10, 4) "19.42"
$ws.Cells.Item(10, 5).Value = "  -0.90%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.0840"
$ws.Cells.Item(11, 5).Value = "  -0.25%  "

Set-TextValue $ws.Cells.Item(12, 4) "1.822.82"
$ws.Cells.Item(12, 5).Value = "  -0.01%  "

Set-TextValue $ws.Cells.Item(13, 4) "1.598.43"
$ws.Cells.Item(13, 5).Value = "  +0.05%  "

Set-TextValue $ws.Cells.Item(14, 4) "4.02"
$ws.Cells.Item(14, 5).Value = "  +0.10%  "

Set-TextValue $ws.Cells.Item(15, 4) "0.523"
$ws.Cells.Item(15, 5).Value = "  +0.36%  "

Set-TextValue $ws.Cells.Item(16, 4) "65.07"
$ws.Cells.Item(16, 5).Value = "  +0.19%  "

Set-TextValue $ws.Cells.Item(17, 4) "26.665.16"

Set-TextValue $ws.Cells.Item(18, 4) "0.0₃0748"
$ws.Cells.Item(18, 5).Value = "  +2.82%  "

$ws.Cells.Item(19, 2).Value = "Dai"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Cells.Item(19, 4) "1.00"
$ws.Cells.Item(19, 5).Value = "  +0.33%  "

$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(20, 4) "209.06"
$ws.Cells.Item(20, 5).Value = "  +0.32%  "

Set-TextValue $ws.Cells.Item(21, 4) "6.98"
$ws.Cells.Item(21, 5).Value = "  +2.74%  "

$ws.Cells.Item(22, 5).Value = "  +0.66%  "

$ws.Cells.Item(23, 5).Value = "  -0.02%  "

Set-TextValue $ws.Cells.Item(24, 4) "8.95"
$ws.Cells.Item(24, 5).Value = "  +0.96%  "

Set-TextValue $ws.Cells.Item(25, 4) "142.91"
$ws.Cells.Item(25, 5).Value = "  -2.07%  "

$ws.Cells.Item(26, 5).Value = "  +0.20%  "

Set-TextValue $ws.Cells.Item(27, 4) "7.10"
$ws.Cells.Item(27, 5).Value = "  -1.68%  "

$ws.Cells.Item(28, 5).Value = "  -1.03%  "

Set-TextValue $ws.Cells.Item(29, 4) "15.30"
$ws.Cells.Item(29, 5).Value = "  +0.26%  "

Set-TextValue $ws.Cells.Item(30, 4) "0.0514"
$ws.Cells.Item(30, 5).Value = "  +1.64%  "

$ws.Cells.Item(31, 5).Value = "  -0.13%  "

Set-TextValue $ws.Cells.Item(32, 4) "3.24"
$ws.Cells.Item(32, 5).Value = "  +0.39%  "

$ws.Cells.Item(33, 5).Value = "  +0.32%  "

Set-TextValue $ws.Cells.Item(34, 4) "1.286.18"
$ws.Cells.Item(34, 5).Value = "  -0.23%  "

Set-TextValue $ws.Cells.Item(35, 4) "0.615"
$ws.Cells.Item(35, 5).Value = "  -6.91%  "

$ws.Cells.Item(36, 5).Value = "  -0.30%  "

Set-TextValue $ws.Cells.Item(37, 4) "1.48"
$ws.Cells.Item(37, 5).Value = "  -0.56%  "

$ws.Cells.Item(38, 5).Value = "  -0.32%  "

$ws.Cells.Item(39, 5).Value = "  +17.09%  "

Set-TextValue $ws.Cells.Item(40, 4) "0.825"
$ws.Cells.Item(40, 5).Value = "  -2.01%  "

Set-TextValue $ws.Cells.Item(41, 4) "5.43"
$ws.Cells.Item(41, 5).Value = "  +0.09%  "

$ws.Cells.Item(42, 5).Value = "  -0.59%  "

$ws.Cells.Item(43, 5).Value = "  -0.57%  "

Set-TextValue $ws.Cells.Item(44, 4) "62.98"
$ws.Cells.Item(44, 5).Value = "  -0.88%  "

Set-TextValue $ws.Cells.Item(45, 4) "1.734.90"
$ws.Cells.Item(45, 5).Value = "  -0.03%  "

Set-TextValue $ws.Cells.Item(46, 4) "90.76"
$ws.Cells.Item(46, 5).Value = "  +0.85%  "

Set-TextValue $ws.Cells.Item(47, 4) "1.56"
$ws.Cells.Item(47, 5).Value = "  -2.79%  "

$ws.Cells.Item(48, 5).Value = "  -0.60%  "

$ws.Cells.Item(49, 5).Value = "  +0.70%  "

$ws.Cells.Item(50, 5).Value = "  +0.28%  "

Set-TextValue $ws.Cells.Item(51, 4) "7.34"
$ws.Cells.Item(51, 5).Value = "  -2.03%  "

